$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text (they are dotted/
# thousand-grouped display strings, not numbers), then update values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.196.11"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.063.25"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.90"
$ws.Range("E5").Value = "  -3.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.02"
$ws.Range("E8").Value = "  +4.76%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.01"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.371.50"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.44"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.10"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.772"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.23"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.059.58"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.345.38"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  +19.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.27"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0808"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "223.42"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.68"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.82"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.130"
$ws.Range("E29").Value = "  +5.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.20"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("E31").Value = "  +6.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.45"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0620"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  +6.59%  "
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.93"
$ws.Range("E38").Value = "  +13.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.32"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.45"
$ws.Range("E42").Value = "  +23.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0960"
$ws.Range("E43").Value = "  +8.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.470.53"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.25"
$ws.Range("E45").Value = "  +6.95%  "
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.03"
$ws.Range("E48").Value = "  +4.61%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("E50").Value = "  +7.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.92"
$ws.Range("E51").Value = "  +1.70%  "
